$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 197.14285
$ws.Range("J9").Value = 145
$ws.Range("L9").Value = 145
$ws.Range("N9").Value = -483

$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = 0

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = 0

$ws.Range("H129").Value = 1613.6875
$ws.Range("I129").Value = 1385.4
$ws.Range("K129").Value = 4156.200000000001
$ws.Range("M129").Value = 843.7999999999993

$ws.Range("H138").Value = 1813
$ws.Range("J138").Value = 2529.5
$ws.Range("L138").Value = 7588.5
$ws.Range("N138").Value = -17868.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4039.8462
$ws.Range("J2").Value = 4220.5
$ws.Range("L2").Value = 4220.5
$ws.Range("N2").Value = -4446.5

$ws.Range("H32").Value = 47700.89
$ws.Range("J32").Value = 155698.58
$ws.Range("L32").Value = 155698.58
$ws.Range("N32").Value = -156272.58

$ws.Range("H45").Value = 633266
$ws.Range("I45").Value = 1123757.5
$ws.Range("J45").Value = 2634.1428
$ws.Range("K45").Value = 1123757.5
$ws.Range("L45").Value = 2634.1428
$ws.Range("M45").Value = -1123380.5
$ws.Range("N45").Value = -3388.1428

$ws.Range("H61").Value = 676.8
$ws.Range("I61").Value = 676.8
$ws.Range("K61").Value = 676.8
$ws.Range("M61").Value = -464.8

$ws.Range("H102").Value = 2098.2
$ws.Range("I102").Value = 2216.0833
$ws.Range("J102").Value = 1921.375
$ws.Range("K102").Value = 2216.0833
$ws.Range("L102").Value = 1921.375
$ws.Range("M102").Value = -594.0832999999998
$ws.Range("N102").Value = -5165.375

$ws.Range("H110").Value = 2512.862
$ws.Range("I110").Value = 2541.8928
$ws.Range("K110").Value = 2541.8928
$ws.Range("M110").Value = -496.8928000000001

$ws.Range("H116").Value = 4039.8462
$ws.Range("J116").Value = 4220.5
$ws.Range("L116").Value = 4220.5
$ws.Range("N116").Value = -8808.5

$ws.Range("H132").Value = 11028.619
$ws.Range("I132").Value = 12261.223
$ws.Range("J132").Value = 3633
$ws.Range("K132").Value = 36783.669
$ws.Range("L132").Value = 10899
$ws.Range("M132").Value = -34253.669
$ws.Range("N132").Value = -15959

$ws.Range("H136").Value = 676.8
$ws.Range("I136").Value = 676.8
$ws.Range("K136").Value = 2030.4
$ws.Range("M136").Value = 519.6000000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4039.8462
$ws.Range("J3").Value = 4220.5
$ws.Range("L3").Value = 4220.5
$ws.Range("N3").Value = -4448.5

$ws.Range("H94").Value = 588.75
$ws.Range("I94").Value = 567.9
$ws.Range("K94").Value = 567.9
$ws.Range("M94").Value = -116.9

$ws.Range("H105").Value = 5646.222
$ws.Range("I105").Value = 6892.5
$ws.Range("J105").Value = 3153.6667
$ws.Range("K105").Value = 6892.5
$ws.Range("L105").Value = 3153.6667
$ws.Range("M105").Value = -5145.5
$ws.Range("N105").Value = -6647.6667

$ws.Range("H107").Value = 22770.723
$ws.Range("I107").Value = 30956.383
$ws.Range("K107").Value = 30956.383
$ws.Range("M107").Value = -29036.383

$ws.Range("H134").Value = 1486.9524
$ws.Range("I134").Value = 1484.8334
$ws.Range("K134").Value = 4454.5002
$ws.Range("M134").Value = -1919.5002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 12398.923
$ws.Range("I99").Value = 20098
$ws.Range("J99").Value = 3416.6667
$ws.Range("K99").Value = 20098
$ws.Range("L99").Value = 3416.6667
$ws.Range("M99").Value = -18600
$ws.Range("N99").Value = -6412.6667

$ws.Range("H105").Value = 1362.5
$ws.Range("I105").Value = 1166.6666
$ws.Range("K105").Value = 1166.6666
$ws.Range("M105").Value = 580.3334

$ws.Range("H107").Value = 2429.6365
$ws.Range("J107").Value = 2236.4
$ws.Range("L107").Value = 2236.4
$ws.Range("N107").Value = -6076.4

$ws.Range("H126").Value = 12398.923
$ws.Range("I126").Value = 20098
$ws.Range("J126").Value = 3416.6667
$ws.Range("K126").Value = 60294
$ws.Range("L126").Value = 10250.0001
$ws.Range("M126").Value = -57824
$ws.Range("N126").Value = -15190.0001

$ws.Range("H132").Value = 6462.5
$ws.Range("I132").Value = 6462.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 19387.5
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -16857.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 183
$ws.Range("I23").Value = 87.5
$ws.Range("J23").Value = 246.66667
$ws.Range("K23").Value = 262.5
$ws.Range("L23").Value = 740.00001
$ws.Range("M23").Value = -27.5
$ws.Range("N23").Value = -1210.00001

$ws.Range("H107").Value = 713.65
$ws.Range("J107").Value = 1004.2222
$ws.Range("L107").Value = 3012.6666
$ws.Range("N107").Value = -6852.6666

$ws.Range("H113").Value = 542.4545000000001
$ws.Range("J113").Value = 496.33334
$ws.Range("L113").Value = 1489.00002
$ws.Range("N113").Value = -5829.000019999999

$ws.Range("H131").Value = 1519.1538
$ws.Range("I131").Value = 900
$ws.Range("J131").Value = 1570.75
$ws.Range("K131").Value = 2700
$ws.Range("L131").Value = 4712.25
$ws.Range("M131").Value = 2340
$ws.Range("N131").Value = -14792.25

$ws.Range("H132").Value = 1104.0667
$ws.Range("I132").Value = 961.3333
$ws.Range("K132").Value = 8651.9997
$ws.Range("M132").Value = -6121.9997

$ws.Range("H138").Value = 2129
$ws.Range("I138").Value = 1188
$ws.Range("K138").Value = 3564
$ws.Range("M138").Value = 1576

$ws.Range("H140").Value = 2374.8096
$ws.Range("I140").Value = 1993.55
$ws.Range("K140").Value = 5980.65
$ws.Range("M140").Value = -800.6499999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 6631.6665
$ws.Range("I10").Value = 5650
$ws.Range("K10").Value = 5650
$ws.Range("M10").Value = -5481

$ws.Range("H80").Value = 5241.222
$ws.Range("I80").Value = 3896
$ws.Range("J80").Value = 6922.75
$ws.Range("K80").Value = 3896
$ws.Range("L80").Value = 6922.75
$ws.Range("M80").Value = -2898
$ws.Range("N80").Value = -8918.75

$ws.Range("H83").Value = 5241.222
$ws.Range("I83").Value = 3896
$ws.Range("J83").Value = 6922.75
$ws.Range("K83").Value = 19480
$ws.Range("L83").Value = 34613.75
$ws.Range("M83").Value = -14488
$ws.Range("N83").Value = -44597.75

$ws.Range("H88").Value = 29717
$ws.Range("J88").Value = 29717
$ws.Range("L88").Value = 29717
$ws.Range("N88").Value = -30619

$ws.Range("H91").Value = 29717
$ws.Range("J91").Value = 29717
$ws.Range("L91").Value = 29717
$ws.Range("N91").Value = -32837

$ws.Range("H122").Value = 152734.88
$ws.Range("I122").Value = 173473.66
$ws.Range("K122").Value = 520420.98
$ws.Range("M122").Value = -517970.98

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1949.6666
$ws.Range("J22").Value = 1949.6666
$ws.Range("L22").Value = 1949.6666
$ws.Range("N22").Value = -2539.6666

$ws.Range("H27").Value = 1949.6666
$ws.Range("J27").Value = 1949.6666
$ws.Range("L27").Value = 1949.6666
$ws.Range("N27").Value = -2163.6666

$ws.Range("H82").Value = 3860.5
$ws.Range("I82").Value = 3721
$ws.Range("J82").Value = 4000
$ws.Range("K82").Value = 3721
$ws.Range("L82").Value = 4000
$ws.Range("M82").Value = -3360
$ws.Range("N82").Value = -4722

$ws.Range("H85").Value = 3860.5
$ws.Range("I85").Value = 3721
$ws.Range("J85").Value = 4000
$ws.Range("K85").Value = 3721
$ws.Range("L85").Value = 4000
$ws.Range("M85").Value = -2473
$ws.Range("N85").Value = -6496

$ws.Range("H93").Value = 39180.223
$ws.Range("I93").Value = 2455.5715
$ws.Range("K93").Value = 2455.5715
$ws.Range("M93").Value = -1207.5715

$ws.Range("H136").Value = 2126.2727
$ws.Range("I136").Value = 1515.5834
$ws.Range("J136").Value = 4874.375
$ws.Range("K136").Value = 4546.7502
$ws.Range("L136").Value = 14623.125
$ws.Range("M136").Value = -1996.7502
$ws.Range("N136").Value = -19723.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1508001
$ws.Range("I14").Value = 3000250
$ws.Range("K14").Value = 3000250
$ws.Range("M14").Value = -3000082

$ws.Range("H62").Value = 8621.357
$ws.Range("I62").Value = 8558.5
$ws.Range("K62").Value = 8558.5
$ws.Range("M62").Value = -7934.5

$ws.Range("H65").Value = 8621.357
$ws.Range("I65").Value = 8558.5
$ws.Range("K65").Value = 42792.5
$ws.Range("M65").Value = -39672.5

$ws.Range("H107").Value = 27778516
$ws.Range("I107").Value = 688
$ws.Range("K107").Value = 2064
$ws.Range("M107").Value = -144

$ws.Range("H136").Value = 639.0571
$ws.Range("I136").Value = 639.0571
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 1917.1713
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = 632.8287
